# fixed score addup rounding
# The "total score average" columns (E, J, O) were previously stored as
# rounded integers. They should instead be the precise sum of the two
# "part average" columns that precede them:
#   E = C + D   (191215)
#   J = H + I   (200112)
#   O = M + N   (200530)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 36 }

for ($r = 2; $r -le $lastRow; $r++) {
    $c = $ws.Cells.Item($r, 3).Value2   # C
    $d = $ws.Cells.Item($r, 4).Value2   # D
    $h = $ws.Cells.Item($r, 8).Value2   # H
    $i = $ws.Cells.Item($r, 9).Value2   # I
    $m = $ws.Cells.Item($r, 13).Value2  # M
    $n = $ws.Cells.Item($r, 14).Value2  # N

    if ($null -eq $c) { $c = 0 }
    if ($null -eq $d) { $d = 0 }
    if ($null -eq $h) { $h = 0 }
    if ($null -eq $i) { $i = 0 }
    if ($null -eq $m) { $m = 0 }
    if ($null -eq $n) { $n = 0 }

    $ws.Cells.Item($r, 5).Value2 = $c + $d    # E
    $ws.Cells.Item($r, 10).Value2 = $h + $i   # J
    $ws.Cells.Item($r, 15).Value2 = $m + $n   # O
}
